$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds the IG metadata table (Property / Value columns)
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 is the "Name" property; its Value cell (B4) was empty and now gets
# the generated FHIR resource name.
$ws.Range("B4").Value = "StatutprofessionnelssaVs"

# Row 8 is the "Date" property; update the generation timestamp.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
